# x1469 Humans have been renamed "Homo sapiens (Human)"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K3").Value = "Homo sapiens (Human)"
$ws.Range("K4").Value = "Homo sapiens (Human)"
$ws.Range("K5").Value = "Homo sapiens (Human)"

# The longer species label now wraps across more lines in the (narrow,
# wrap-text) Species column, so the rows grow taller to fit it.
$ws.Rows("3").RowHeight = 52
$ws.Rows("4").RowHeight = 52
$ws.Rows("5").RowHeight = 52

$ws.Range("K5").Select()
